# "add nodejs http agent"
#
# The deck has a slide (index 10) that reads "Node.js Advanced" / "RPC"
# (a gradient rectangle group with two text boxes). The edit inserts a new,
# similar slide right before it that reads "Node.js Advanced" / "http.Agent"
# with a different accent gradient (teal -> blue instead of pink -> purple).
# All following slides shift down by one position.
#
# Strategy: duplicate slide 10 (the duplicate lands immediately after the
# original, i.e. at index 11) and then re-paint the *original* slide 10 in
# place with the new gradient colors and the new "http.Agent" caption, so the
# final order is: ... slide9, [new http.Agent slide], [old RPC slide], ...

$p = $ppt.ActivePresentation

$source = $p.Slides.Item(10)
$source.Duplicate() | Out-Null

$newSlide = $p.Slides.Item(10)
$grp = $newSlide.Shapes.Item(1)

$rect = $grp.GroupItems.Item(1)
$titleBox = $grp.GroupItems.Item(2)
$subtitleBox = $grp.GroupItems.Item(3)

# Gradient: FF7AF5/513162 -> 77EFD8/45BAF2
$rect.Fill.GradientStops.Item(1).Color.RGB = 14217079
$rect.Fill.GradientStops.Item(2).Color.RGB = 15907397

# Title text stays "Node.js Advanced" (unchanged wording).
$titleBox.TextFrame.TextRange.Text = "Node.js Advanced"

# Subtitle: "RPC" -> "http.Agent", color 000080 -> FEB692
$subtitleBox.TextFrame.TextRange.Text = "http.Agent"
$subtitleBox.TextFrame.TextRange.Font.Color.RGB = 9615102
